# Add "Brake or Coast?" data to the PDP sheet (column K), and make PDP the
# active/selected tab (it was previously PCM).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDP")

# Header
$ws.Range("K1").Value = "Brake or Coast?"

# Per-device values. Most motor controllers are wired for brake mode; the
# two climb Victor SPXs (rows 10-11) are still undecided.
$ws.Range("K2").Value  = "Brake"
$ws.Range("K3").Value  = "Brake"
$ws.Range("K4").Value  = "Brake"
$ws.Range("K5").Value  = "Brake"
$ws.Range("K6").Value  = "Brake"
$ws.Range("K7").Value  = "Brake"
$ws.Range("K8").Value  = "Brake"
$ws.Range("K9").Value  = "Brake"
$ws.Range("K10").Value = "???"
$ws.Range("K11").Value = "???"
$ws.Range("K12").Value = "Brake"
$ws.Range("K13").Value = "Brake"

# Match the column width that Excel's bestFit would have produced for the
# new, narrower column J and the new column K (values chosen so the
# engine's internal width quantization lands on the closest stored width).
$ws.Columns("J").ColumnWidth = 15.0
$ws.Columns("K").ColumnWidth = 13.83

# Make PDP the active sheet/tab, with K14 as the selected cell (this also
# clears tabSelected on whichever sheet was previously active).
$ws.Activate()
$ws.Range("K14").Select()
